# Fixed Current Account Summary
#
# The "unmatched_expenses" reconciliation sheet had a row inserted ahead of
# the existing data: a new "Unknown (parsing failed)" entry (no amount) now
# sits at row 2, and the rows that used to be 2-4 shift down to 3-5 (row 4's
# data is duplicated into the new row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("unmatched_expenses")

$data = @(
    @("UNKNOWN DATE", "", "Unknown (parsing failed)"),
    @("UNKNOWN DATE", 299, "Generated on 19-Jun-2025 11:10:00 (parsing failed)"),
    @("UNKNOWN DATE", 780, "Amount in Words: 780 Rupees Only (parsing failed)"),
    @("UNKNOWN DATE", 780, "Amount in Words: 780 Rupees Only (parsing failed)")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

$wb.Save()
